$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("karamjeet")

# --- header text first, so new shared strings are interned in "id, category,
#     start_date, end_date" order (matches the source sheet's si order) ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "start_date"
$ws.Range("E1").Value = "end_date"

# --- data values ---
$ws.Range("A2").Value = 146
$ws.Range("B2").Value = "Shankar1"
$ws.Range("C2").Value = "null"

# --- style 5: plain centered, applied to the whole data block first ---
$ws.Range("A2:E8").HorizontalAlignment = -4108

# --- style 6: date centered (numFmtId 14 + center), built on D2 then copied to E2 ---
$ws.Range("D2").NumberFormatLocal = "mm-dd-yy"
$ws.Range("D2").Value = Get-Date -Year 2022 -Month 4 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Value = Get-Date -Year 2022 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0
$excel.CutCopyMode = 0

# --- style 7: bold centered header ---
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").Font.Bold = $true

# --- column widths ---
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 14.36328125
$ws.Columns.Item(5).ColumnWidth = 17.26953125

# --- selection / active sheet ---
$ws.Activate()
$ws.Range("B5").Select()
